# Adds the new custom styles (ContactInfo, SkillCategory, SkillItems,
# SkillHighlight, SkillLevel) used by the "minimal" CV template to the
# document's style sheet.

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1, wdStyleTypeCharacter = 2
# wdLineSpaceMultiple = 5
# Colors are Word OLE COLORREF (0x00BBGGRR), computed from the target
# hex RGB: 1F2937 -> 3615007, 6B7280 -> 8417899

# --- Contact Info (paragraph) ---
$s = $d.Styles.Add("Contact Info", 1)
$s.ParagraphFormat.SpaceBefore = 0
$s.ParagraphFormat.SpaceAfter = 5.1
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 18
$s.Font.Name = "Liberation Sans"
$s.Font.Bold = $false
$s.Font.Color = 3615007
$s.Font.Size = 9

# --- Skill Category (paragraph) ---
$s = $d.Styles.Add("Skill Category", 1)
$s.ParagraphFormat.SpaceBefore = 0
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 14.4
$s.Font.Name = "Liberation Sans"
$s.Font.Bold = $true
$s.Font.Color = 8417899
$s.Font.Size = 10

# --- Skill Items (paragraph) ---
$s = $d.Styles.Add("Skill Items", 1)
$s.ParagraphFormat.SpaceBefore = 0
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 14.4
$s.Font.Name = "Liberation Sans"
$s.Font.Bold = $false
$s.Font.Color = 3615007
$s.Font.Size = 10

# --- Skill Highlight (character) ---
$s = $d.Styles.Add("Skill Highlight", 2)
$s.Font.Name = "Liberation Sans"
$s.Font.Bold = $true
$s.Font.Color = 3615007
$s.Font.Size = 10

# --- Skill Level (character) ---
$s = $d.Styles.Add("Skill Level", 2)
$s.Font.Name = "Liberation Sans"
$s.Font.Bold = $false
$s.Font.Color = 8417899
$s.Font.Size = 10
